# Daily attendance processing - 2025-11-20 17:45:10
#
# Normalizes the "Recorded By" column (G) so that when the value is a
# comma-separated list that begins with "System", the "System" token is
# moved from the first position to the second position in the list
# (i.e. "System, X, ..." becomes "X, System, ...").
#
# Only the specific rows touched by this processing run are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,10,12,13,14,15,18,19,20,21,22,24,26,28,29,30,31,32,33,34,36,38,39,40,41,44,45,46,47,48,50,52,54,55,56,57,58,59,60,62,64,65,66,67,70,71,72,73,74,76,78,80,81,82,83,84,85,86,90,92,99,101,106,107,108,109,110,111,112,116,118,125,127,132,133,134,135,136,137,138,142,144,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -and $value.StartsWith("System, ")) {
        $parts = $value.Split(", ")
        if ($parts.Count -gt 1 -and $parts[0] -eq "System") {
            $newParts = New-Object System.Collections.ArrayList
            [void]$newParts.Add($parts[1])
            [void]$newParts.Add("System")
            for ($i = 2; $i -lt $parts.Count; $i++) {
                [void]$newParts.Add($parts[$i])
            }
            $cell.Value2 = ($newParts -join ", ")
        }
    }
}
